$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume snapshot data (Coin, Link, Price, Volume(1h)).
# Column B:E are forced to text format ("@") before assignment so that
# numeric-looking strings (e.g. "234.07", "30.157.55") are preserved
# exactly as text rather than being auto-converted to numbers by Excel.
$rows = @(
  @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.157.55', '  -0.56%  '),
  @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.861.88', '  -0.45%  '),
  @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  -0.04%  '),
  @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '233.91', '  -0.91%  '),
  @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.001', '  +0.02%  '),
  @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4678', '  -0.50%  '),
  @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2850', '  -1.27%  '),
  @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06472', '  -2.21%  '),
  @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '21.01', '  -3.11%  '),
  @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07751', '  -3.62%  '),
  @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.878.37', '  +0.40%  '),
  @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '93.47', '  -4.12%  '),
  @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6810', '  -0.83%  '),
  @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.041', '  -1.93%  '),
  @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '267.28', '  -1.70%  '),
  @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.147.88', '  -0.55%  '),
  @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.31', '  -5.77%  '),
  @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007602', '  -1.76%  '),
  @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.001', '  +0.01%  '),
  @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.121.30', '  +0.19%  '),
  @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  -0.06%  '),
  @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.133', '  -3.43%  '),
  @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.097', '  -1.98%  '),
  @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.340', '  -0.77%  '),
  @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '164.94', '  -2.17%  '),
  @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.49', '  -2.48%  '),
  @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.886', '  -3.61%  '),
  @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09947', '  +0.52%  '),
  @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.363', '  -0.79%  '),
  @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.449', '  -1.34%  '),
  @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.215', '  -3.70%  '),
  @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.992', '  -2.30%  '),
  @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04667', '  -0.93%  '),
  @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.114', '  -1.67%  '),
  @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6891', '  -1.87%  '),
  @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.715', '  +0.38%  '),
  @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01833', '  -2.79%  '),
  @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.761', '  +4.08%  '),
  @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.301', '  -0.19%  '),
  @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '71.18', '  -2.30%  '),
  @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  +0.00%  '),
  @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.883', '  -3.92%  '),
  @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8303', '  -1.67%  '),
  @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '102.04', '  -1.28%  '),
  @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4045', '  -3.07%  '),
  @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.127', '  -1.53%  '),
  @('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '931.24', '  +0.22%  '),
  @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.940', '  -2.17%  '),
  @('Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '34.01', '  -1.40%  '),
  @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05572', '  -2.00%  ')
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $rowVals = $rows[$i]
  $rng = $ws.Range("B" + $r + ":E" + $r)
  $rng.NumberFormat = "@"
  $ws.Cells.Item($r, 2).Value = $rowVals[0]
  $ws.Cells.Item($r, 3).Value = $rowVals[1]
  $ws.Cells.Item($r, 4).Value = $rowVals[2]
  $ws.Cells.Item($r, 5).Value = $rowVals[3]
}

Write-Host "Updated cryptos list"
